$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 305
$ws.Range("C2").Value = "2022-11-13 15:57"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "'4673"
$ws.Range("G2").Value = "2022-11-13 16:04"

# Row 3
$ws.Range("B3").Value = 287
$ws.Range("C3").Value = "2022-11-13 15:57"
$ws.Range("E3").Value = "'3616"
$ws.Range("G3").Value = "2022-11-13 16:04"
